$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column indices corresponding to B,C,D,F,G,H,L,M,N,O (E,I,J,K unchanged)
$colIndices = @(2,3,4,6,7,8,12,13,14,15)

# Data rows for sheet rows 2..25, in the same column order as $colIndices
$rows = @(
    @(0.7837220181217504,0.4100876787769749,0.02096903008951045,0.5383616254118238,0.3800234734114838,0.5335069827745826,0.2850609495427392,0.2086486547752315,1.220696161227117,1.768812778141466),
    @(0.7042054461034581,0.4045696761881459,0.01891952754240123,0.5341870260095263,0.3772373196669037,0.5358212346289974,0.2831093610163222,0.195719874688919,1.229440667639736,1.767486867966284),
    @(0.6554113477591272,0.4011981938140821,0.01765259897646132,0.5319974169819233,0.3758316283522234,0.5375094016707465,0.282060902933857,0.1878451596558719,1.235274187941556,1.767869265150949),
    @(0.6355361125208958,0.3998286419261774,0.01713419994661081,0.5311990783837643,0.375335407922698,0.5382645747412838,0.2816713804305451,0.184652361704515,1.237768291165871,1.76832585561435),
    @(0.6322364138470959,0.3996014972417328,0.01704799325784023,0.5310721888197278,0.3752576353463013,0.5383940324356402,0.2816089811815843,0.1841231856297725,1.238189499239652,1.76841983402403),
    @(0.6551432661263448,0.401179705684342,0.0176456161934837,0.5319862699336753,0.3758246260768914,0.5375193139290388,0.2820554968416786,0.1878020345343359,1.235307350806821,1.767874205272079),
    @(0.7562993955312152,0.4081817346651491,0.02026414766629614,0.5368446758400012,0.3789994455040642,0.5342495012395148,0.2843569790059348,0.2041777357302408,1.223614966348975,1.768107195345408),
    @(0.9548492758156044,0.4220370671693843,0.02533040209442561,0.5493379430896255,0.3876504258175828,0.5299561795321779,0.2900570165234342,0.236787287275547,1.204365228007397,1.778065811520491),
    @(1.100782146520032,0.4322836473204887,0.02900965444706571,0.5603289341288971,0.395493037149933,0.5280921861893972,0.2949665998233399,0.2610403824234879,1.192458209546999,1.791190072170309),
    @(1.167172992667645,0.4369578581387827,0.03067391433651068,0.5657234977607501,0.3993855223134233,0.5275241394683121,0.2973565009629908,0.2721362146557027,1.187525475731732,1.798425206377743),
    @(1.192313031413221,0.4387295570904257,0.03130274193423332,0.5678230723281672,0.4009063347502888,0.5273492547718917,0.298283954745898,0.2763467858568731,1.185727048169255,1.80134704282537),
    @(1.186898730587814,0.4383479184835437,0.03116737506656619,0.5673683663272584,0.4005767169813765,0.527385130840571,0.2980832133913793,0.2754395754792753,1.186111282332206,1.800709673836963),
    @(1.169241302016815,0.4371035843656728,0.03072567643567226,0.5658950931980229,0.3995097016977951,0.5275089456342528,0.2974323535291603,0.2724824454603052,1.187376125755044,1.798661938171762),
    @(1.158425477434378,0.43634160649097,0.03045494108906865,0.5650000643752975,0.3988622233996466,0.5275900231223005,0.2970366050848838,0.2706722596095901,1.18815992636307,1.797431353832224),
    @(1.096443326948588,0.4319784216827145,0.02890069820275443,0.5599843302876906,0.395245199405025,0.5281349393239196,0.2948135587527929,0.2603164908107587,1.192790304051002,1.790742705220566),
    @(1.058419570408773,0.4293049399186941,0.02794477691837471,0.557008442878967,0.3931095444851138,0.5285408966574039,0.293489837687261,0.2539795244253753,1.195754751153792,1.786963511304748),
    @(1.036549894537529,0.42776845799456,0.02739406870709615,0.5553339422025942,0.3919117419944484,0.5288007390153098,0.2927432028604926,0.2503406108137014,1.197505367634385,1.784908862213882),
    @(1.029145343616278,0.4272484494064912,0.02720745691886606,0.5547733653785798,0.3915114339319103,0.528893242990776,0.2924929374894276,0.249109564890496,1.198105921248093,1.784233634035189),
    @(1.06246721463333,0.4295894102829863,0.02804662846609318,0.5573213860864641,0.3933337240081727,0.5284949552362406,0.2936292254882318,0.2546534920640156,1.195434467578615,1.787353491203078),
    @(1.174427748586368,0.4374690314731993,0.03085545205586726,0.566326288409293,0.3998218386701922,0.5274714868157844,0.2976229180020056,0.2733507887263329,1.187002725234414,1.799258465742099),
    @(1.247595462338438,0.4426285046534986,0.03268305473956445,0.5725424182778625,0.4043351030223619,0.5270370250040628,0.300363857305328,0.2856218552187642,1.181897109328908,1.808100244633607),
    @(1.208545422669999,0.4398739734966455,0.03170838308503932,0.5691944713713468,0.4019012850455539,0.5272474633717508,0.2988890138125129,0.279067939843479,1.184585039934362,1.803284068290594),
    @(1.060637303571355,0.42946079954136,0.02800058491856561,0.5571797911533238,0.3932322789628699,0.5285156429537636,0.2935661634366085,0.2543487779054487,1.195579123503776,1.787176813555675),
    @(0.9011221946013279,0.4182764252749109,0.02396730743405584,0.5456402320506655,0.3850496132110948,0.5308909470642931,0.2883881276448861,0.2279131790217477,1.209179674717191,1.774352937817866)
)

$startRow = 2
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowVals = $rows[$r]
    $sheetRow = $startRow + $r
    for ($i = 0; $i -lt $colIndices.Length; $i++) {
        $ws.Cells.Item($sheetRow, $colIndices[$i]).Value = $rowVals[$i]
    }
}

Write-Host "case with 380 kV done"